$wb = $excel.ActiveWorkbook

# --- Rename the "Include from Medication statu" worksheet ---
$wsInclude = $wb.Worksheets.Item("Include from Medication statu")
$wsInclude.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value
$ws.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"

# Insert a new row for "Jurisdiction" right after "Contact" (row 10), before "Description" (row 11)
$ws.Rows.Item(11).Insert()

# Copy formatting from the row above (Contact) so the new row matches the sheet's styling
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" row
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
